$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H12").Value = 1
$ws.Range("H13").Value = 0
$ws.Range("H20").Value = 1
$ws.Range("H21").Value = 0
